# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets now that the
# handback xliffs have come back in sync with en-US, flips the Overview
# status text, and widens the columns that now hold the longer status /
# filename strings.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: flip the per-locale status cells --------------------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: record the target/handback file + datetime -------------
$wsZhCn.Range("I2").Value = "2fa36fe8-2159-4d3a-b069-84e3eee85dcd.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/1460fa40ee0ed3f6c1c370cf10a52b9c282209ed/e2e/2fa36fe8-2159-4d3a-b069-84e3eee85dcd.md", "", "", "2fa36fe8-2159-4d3a-b069-84e3eee85dcd.md")
$wsZhCn.Range("J2").Value = "2fa36fe8-2159-4d3a-b069-84e3eee85dcd.b53354ada97b04e6c02a94ba6c4b7b91e3668759.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-14 01:19:14"

# --- de-de sheet: record the target/handback file + datetime -------------
$wsDeDe.Range("I2").Value = "2fa36fe8-2159-4d3a-b069-84e3eee85dcd.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/1460fa40ee0ed3f6c1c370cf10a52b9c282209ed/e2e/2fa36fe8-2159-4d3a-b069-84e3eee85dcd.md", "", "", "2fa36fe8-2159-4d3a-b069-84e3eee85dcd.md")
$wsDeDe.Range("J2").Value = "2fa36fe8-2159-4d3a-b069-84e3eee85dcd.b53354ada97b04e6c02a94ba6c4b7b91e3668759.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-14 01:19:24"

# --- Column widths: the Status/Target/Handback columns grew to fit the
#     longer "Handed back: in sync with en-US" / *.md / *.xlf strings -----
$wsOverview.Range("E1").ColumnWidth = 29.166666666666668
$wsOverview.Range("F1").ColumnWidth = 29.166666666666668

$wsZhCn.Range("C1").ColumnWidth = 29.166666666666668
$wsZhCn.Range("I1").ColumnWidth = 39.166666666666664
$wsZhCn.Range("J1").ColumnWidth = 39.166666666666664

$wsDeDe.Range("C1").ColumnWidth = 29.166666666666668
$wsDeDe.Range("I1").ColumnWidth = 39.166666666666664
$wsDeDe.Range("J1").ColumnWidth = 39.166666666666664
